# Brasil.xlsx data update: append 2020-04-15 to the "Confirmados" and
# "Mortes" daily time series (row 52), one new data point per UF column.

$wb = $excel.ActiveWorkbook

function Set-NewRow {
    param($ws, $rowNum, $dateText, $values)

    # Column A holds the date as plain text (matching the existing rows,
    # which are shared-string "yyyy-mm-dd" text, not real dates). Force
    # text entry with a temporary Text number format, then drop the
    # formatting again so the new cell ends up styleless, like its peers.
    $dateCell = $ws.Range("A" + $rowNum)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dateText
    $dateCell.ClearFormats()

    foreach ($col in $values.Keys) {
        $ws.Range($col + $rowNum).Value = $values[$col]
    }
}

$confirmados = $wb.Worksheets.Item("Confirmados")
$mortes = $wb.Worksheets.Item("Mortes")

# New row for "Confirmados" (sheet1), columns B:AB = Acre .. Tocantins.
$rowConfirmados = @{
    "B"  = 101
    "C"  = 82
    "D"  = 334
    "E"  = 1554
    "F"  = 807
    "G"  = 2157
    "H"  = 682
    "I"  = 557
    "J"  = 304
    "K"  = 630
    "L"  = 151
    "M"  = 121
    "N"  = 903
    "O"  = 384
    "P"  = 136
    "Q"  = 803
    "R"  = 1484
    "S"  = 75
    "T"  = 3743
    "U"  = 399
    "V"  = 747
    "W"  = 69
    "X"  = 114
    "Y"  = 826
    "Z"  = 11043
    "AA" = 46
    "AB" = 26
}

# New row for "Mortes" (sheet2), columns B:AB = Acre .. Tocantins.
$rowMortes = @{
    "B"  = 3
    "C"  = 5
    "D"  = 7
    "E"  = 106
    "F"  = 27
    "G"  = 116
    "H"  = 17
    "I"  = 18
    "J"  = 15
    "K"  = 34
    "L"  = 4
    "M"  = 4
    "N"  = 30
    "O"  = 21
    "P"  = 16
    "Q"  = 38
    "R"  = 143
    "S"  = 8
    "T"  = 265
    "U"  = 19
    "V"  = 19
    "W"  = 2
    "X"  = 3
    "Y"  = 26
    "Z"  = 778
    "AA" = 4
    "AB" = 0
}

Set-NewRow $confirmados 52 "2020-04-15" $rowConfirmados
Set-NewRow $mortes 52 "2020-04-15" $rowMortes
